$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"235.2675016666667"
$ws.Range("H2").Value = [double]"705.802505"
$ws.Range("I2").Value = [double]"0.5738994362335403"
$ws.Range("J2").Value = [double]"0.5738994362335402"
$ws.Range("M2").Value = [double]"4.925988333333333"
$ws.Range("N2").Value = [double]"14.777965"
$ws.Range("O2").Value = [double]"0.05656988822582037"
$ws.Range("P2").Value = [double]"0.05656988822582035"
$ws.Range("Q2").Value = [double]"1158.924968422481"
$ws.Range("R2").Value = [double]"10430.32471580232"
$ws.Range("S2").Value = [double]"0.0324654269605927"
$ws.Range("T2").Value = [double]"0.03246542696059269"
$ws.Range("G3").Value = [double]"235.2675016666667"
$ws.Range("H3").Value = [double]"705.802505"
$ws.Range("I3").Value = [double]"0.5738994362335403"
$ws.Range("J3").Value = [double]"0.5738994362335402"
$ws.Range("O3").Value = [double]"0.5464678959362861"
$ws.Range("P3").Value = [double]"0.5464678959362861"
$ws.Range("Q3").Value = [double]"11195.27205911631"
$ws.Range("R3").Value = [double]"100757.4485320468"
$ws.Range("S3").Value = [double]"0.3136176173975636"
$ws.Range("T3").Value = [double]"0.3136176173975635"
$ws.Range("G4").Value = [double]"235.2675016666667"
$ws.Range("H4").Value = [double]"705.802505"
$ws.Range("I4").Value = [double]"0.5738994362335403"
$ws.Range("J4").Value = [double]"0.5738994362335402"
$ws.Range("M4").Value = [double]"13.80191933333334"
$ws.Range("N4").Value = [double]"41.40575800000001"
$ws.Range("O4").Value = [double]"0.1585007882996995"
$ws.Range("P4").Value = [double]"0.1585007882996994"
$ws.Range("Q4").Value = [double]"3247.1430797582"
$ws.Range("R4").Value = [double]"29224.28771782379"
$ws.Range("S4").Value = [double]"0.09096351304776924"
$ws.Range("T4").Value = [double]"0.09096351304776921"
$ws.Range("G5").Value = [double]"235.2675016666667"
$ws.Range("H5").Value = [double]"705.802505"
$ws.Range("I5").Value = [double]"0.5738994362335403"
$ws.Range("J5").Value = [double]"0.5738994362335402"
$ws.Range("M5").Value = [double]"12.18796133333333"
$ws.Range("N5").Value = [double]"36.563884"
$ws.Range("O5").Value = [double]"0.1399661476381804"
$ws.Range("P5").Value = [double]"0.1399661476381803"
$ws.Range("Q5").Value = [double]"2867.431213303269"
$ws.Range("R5").Value = [double]"25806.88091972942"
$ws.Range("S5").Value = [double]"0.08032649322133217"
$ws.Range("T5").Value = [double]"0.08032649322133215"
$ws.Range("G6").Value = [double]"235.2675016666667"
$ws.Range("H6").Value = [double]"705.802505"
$ws.Range("I6").Value = [double]"0.5738994362335403"
$ws.Range("J6").Value = [double]"0.5738994362335402"
$ws.Range("M6").Value = [double]"8.576764333333333"
$ws.Range("N6").Value = [double]"25.730293"
$ws.Range("O6").Value = [double]"0.09849527990001386"
$ws.Range("P6").Value = [double]"0.09849527990001385"
$ws.Range("Q6").Value = [double]"2017.833917087107"
$ws.Range("R6").Value = [double]"18160.50525378397"
$ws.Range("S6").Value = [double]"0.05652638560628271"
$ws.Range("T6").Value = [double]"0.0565263856062827"
$ws.Range("I7").Value = [double]"0.3286113026040369"
$ws.Range("J7").Value = [double]"0.3286113026040369"
$ws.Range("M7").Value = [double]"4.925988333333333"
$ws.Range("N7").Value = [double]"14.777965"
$ws.Range("O7").Value = [double]"0.05656988822582037"
$ws.Range("P7").Value = [double]"0.05656988822582035"
$ws.Range("Q7").Value = [double]"663.5933396154754"
$ws.Range("R7").Value = [double]"5972.340056539279"
$ws.Range("S7").Value = [double]"0.0185895046580516"
$ws.Range("T7").Value = [double]"0.0185895046580516"
$ws.Range("I8").Value = [double]"0.3286113026040369"
$ws.Range("J8").Value = [double]"0.3286113026040369"
$ws.Range("O8").Value = [double]"0.5464678959362861"
$ws.Range("P8").Value = [double]"0.5464678959362861"
$ws.Range("S8").Value = [double]"0.1795755271149103"
$ws.Range("T8").Value = [double]"0.1795755271149102"
$ws.Range("I9").Value = [double]"0.3286113026040369"
$ws.Range("J9").Value = [double]"0.3286113026040369"
$ws.Range("M9").Value = [double]"13.80191933333334"
$ws.Range("N9").Value = [double]"41.40575800000001"
$ws.Range("O9").Value = [double]"0.1585007882996995"
$ws.Range("P9").Value = [double]"0.1585007882996994"
$ws.Range("Q9").Value = [double]"1859.294241834393"
$ws.Range("R9").Value = [double]"16733.64817650954"
$ws.Range("S9").Value = [double]"0.05208515050693094"
$ws.Range("T9").Value = [double]"0.05208515050693092"
$ws.Range("I10").Value = [double]"0.3286113026040369"
$ws.Range("J10").Value = [double]"0.3286113026040369"
$ws.Range("M10").Value = [double]"12.18796133333333"
$ws.Range("N10").Value = [double]"36.563884"
$ws.Range("O10").Value = [double]"0.1399661476381804"
$ws.Range("P10").Value = [double]"0.1399661476381803"
$ws.Range("Q10").Value = [double]"1641.873552473081"
$ws.Range("R10").Value = [double]"14776.86197225773"
$ws.Range("S10").Value = [double]"0.04599445809585139"
$ws.Range("T10").Value = [double]"0.04599445809585138"
$ws.Range("I11").Value = [double]"0.3286113026040369"
$ws.Range("J11").Value = [double]"0.3286113026040369"
$ws.Range("M11").Value = [double]"8.576764333333333"
$ws.Range("N11").Value = [double]"25.730293"
$ws.Range("O11").Value = [double]"0.09849527990001386"
$ws.Range("P11").Value = [double]"0.09849527990001385"
$ws.Range("Q11").Value = [double]"1155.399343627806"
$ws.Range("R11").Value = [double]"10398.59409265026"
$ws.Range("S11").Value = [double]"0.03236666222829277"
$ws.Range("T11").Value = [double]"0.03236666222829276"
$ws.Range("G12").Value = [double]"0.325805"
$ws.Range("H12").Value = [double]"0.977415"
$ws.Range("I12").Value = [double]"0.0007947519504286909"
$ws.Range("J12").Value = [double]"0.0007947519504286907"
$ws.Range("M12").Value = [double]"4.925988333333333"
$ws.Range("N12").Value = [double]"14.777965"
$ws.Range("O12").Value = [double]"0.05656988822582037"
$ws.Range("P12").Value = [double]"0.05656988822582035"
$ws.Range("Q12").Value = [double]"1.604911628941667"
$ws.Range("R12").Value = [double]"14.444204660475"
$ws.Range("S12").Value = [double]"4.495902900300378E-05"
$ws.Range("T12").Value = [double]"4.495902900300376E-05"
$ws.Range("G13").Value = [double]"0.325805"
$ws.Range("H13").Value = [double]"0.977415"
$ws.Range("I13").Value = [double]"0.0007947519504286909"
$ws.Range("J13").Value = [double]"0.0007947519504286907"
$ws.Range("O13").Value = [double]"0.5464678959362861"
$ws.Range("P13").Value = [double]"0.5464678959362861"
$ws.Range("Q13").Value = [double]"15.503525082645"
$ws.Range("R13").Value = [double]"139.531725743805"
$ws.Range("S13").Value = [double]"0.0004343064261420263"
$ws.Range("T13").Value = [double]"0.0004343064261420262"
$ws.Range("G14").Value = [double]"0.325805"
$ws.Range("H14").Value = [double]"0.977415"
$ws.Range("I14").Value = [double]"0.0007947519504286909"
$ws.Range("J14").Value = [double]"0.0007947519504286907"
$ws.Range("M14").Value = [double]"13.80191933333334"
$ws.Range("N14").Value = [double]"41.40575800000001"
$ws.Range("O14").Value = [double]"0.1585007882996995"
$ws.Range("P14").Value = [double]"0.1585007882996994"
$ws.Range("Q14").Value = [double]"4.496734328396668"
$ws.Range("R14").Value = [double]"40.47060895557001"
$ws.Range("S14").Value = [double]"0.0001259688106456712"
$ws.Range("T14").Value = [double]"0.0001259688106456711"
$ws.Range("G15").Value = [double]"0.325805"
$ws.Range("H15").Value = [double]"0.977415"
$ws.Range("I15").Value = [double]"0.0007947519504286909"
$ws.Range("J15").Value = [double]"0.0007947519504286907"
$ws.Range("M15").Value = [double]"12.18796133333333"
$ws.Range("N15").Value = [double]"36.563884"
$ws.Range("O15").Value = [double]"0.1399661476381804"
$ws.Range("P15").Value = [double]"0.1399661476381803"
$ws.Range("Q15").Value = [double]"3.970898742206667"
$ws.Range("R15").Value = [double]"35.73808867986001"
$ws.Range("S15").Value = [double]"0.0001112383688294339"
$ws.Range("T15").Value = [double]"0.0001112383688294339"
$ws.Range("G16").Value = [double]"0.325805"
$ws.Range("H16").Value = [double]"0.977415"
$ws.Range("I16").Value = [double]"0.0007947519504286909"
$ws.Range("J16").Value = [double]"0.0007947519504286907"
$ws.Range("M16").Value = [double]"8.576764333333333"
$ws.Range("N16").Value = [double]"25.730293"
$ws.Range("O16").Value = [double]"0.09849527990001386"
$ws.Range("P16").Value = [double]"0.09849527990001385"
$ws.Range("Q16").Value = [double]"2.794352703621667"
$ws.Range("R16").Value = [double]"25.149174332595"
$ws.Range("S16").Value = [double]"7.827931580855586E-05"
$ws.Range("T16").Value = [double]"7.827931580855582E-05"
$ws.Range("G17").Value = [double]"39.46134166666666"
$ws.Range("H17").Value = [double]"118.384025"
$ws.Range("I17").Value = [double]"0.09625996610278018"
$ws.Range("J17").Value = [double]"0.09625996610278018"
$ws.Range("M17").Value = [double]"4.925988333333333"
$ws.Range("N17").Value = [double]"14.777965"
$ws.Range("O17").Value = [double]"0.05656988822582037"
$ws.Range("P17").Value = [double]"0.05656988822582035"
$ws.Range("Q17").Value = [double]"194.3861086676805"
$ws.Range("R17").Value = [double]"1749.474978009125"
$ws.Range("S17").Value = [double]"0.005445415523055532"
$ws.Range("T17").Value = [double]"0.005445415523055531"
$ws.Range("G18").Value = [double]"39.46134166666666"
$ws.Range("H18").Value = [double]"118.384025"
$ws.Range("I18").Value = [double]"0.09625996610278018"
$ws.Range("J18").Value = [double]"0.09625996610278018"
$ws.Range("O18").Value = [double]"0.5464678959362861"
$ws.Range("P18").Value = [double]"0.5464678959362861"
$ws.Range("Q18").Value = [double]"1877.779347536075"
$ws.Range("R18").Value = [double]"16900.01412782467"
$ws.Range("S18").Value = [double]"0.0526029811390845"
$ws.Range("T18").Value = [double]"0.0526029811390845"
$ws.Range("G19").Value = [double]"39.46134166666666"
$ws.Range("H19").Value = [double]"118.384025"
$ws.Range("I19").Value = [double]"0.09625996610278018"
$ws.Range("J19").Value = [double]"0.09625996610278018"
$ws.Range("M19").Value = [double]"13.80191933333334"
$ws.Range("N19").Value = [double]"41.40575800000001"
$ws.Range("O19").Value = [double]"0.1585007882996995"
$ws.Range("P19").Value = [double]"0.1585007882996994"
$ws.Range("Q19").Value = [double]"544.642254468439"
$ws.Range("R19").Value = [double]"4901.78029021595"
$ws.Range("S19").Value = [double]"0.01525728050899301"
$ws.Range("T19").Value = [double]"0.015257280508993"
$ws.Range("G20").Value = [double]"39.46134166666666"
$ws.Range("H20").Value = [double]"118.384025"
$ws.Range("I20").Value = [double]"0.09625996610278018"
$ws.Range("J20").Value = [double]"0.09625996610278018"
$ws.Range("M20").Value = [double]"12.18796133333333"
$ws.Range("N20").Value = [double]"36.563884"
$ws.Range("O20").Value = [double]"0.1399661476381804"
$ws.Range("P20").Value = [double]"0.1399661476381803"
$ws.Range("Q20").Value = [double]"480.9533063947889"
$ws.Range("R20").Value = [double]"4328.5797575531"
$ws.Range("S20").Value = [double]"0.01347313662718797"
$ws.Range("T20").Value = [double]"0.01347313662718796"
$ws.Range("G21").Value = [double]"39.46134166666666"
$ws.Range("H21").Value = [double]"118.384025"
$ws.Range("I21").Value = [double]"0.09625996610278018"
$ws.Range("J21").Value = [double]"0.09625996610278018"
$ws.Range("M21").Value = [double]"8.576764333333333"
$ws.Range("N21").Value = [double]"25.730293"
$ws.Range("O21").Value = [double]"0.09849527990001386"
$ws.Range("P21").Value = [double]"0.09849527990001385"
$ws.Range("Q21").Value = [double]"338.4506277521472"
$ws.Range("R21").Value = [double]"3046.055649769325"
$ws.Range("S21").Value = [double]"0.00948115230445918"
$ws.Range("T21").Value = [double]"0.009481152304459178"
$ws.Range("G22").Value = [double]"0.178139"
$ws.Range("H22").Value = [double]"0.534417"
$ws.Range("I22").Value = [double]"0.0004345431092138444"
$ws.Range("J22").Value = [double]"0.0004345431092138443"
$ws.Range("M22").Value = [double]"4.925988333333333"
$ws.Range("N22").Value = [double]"14.777965"
$ws.Range("O22").Value = [double]"0.05656988822582037"
$ws.Range("P22").Value = [double]"0.05656988822582035"
$ws.Range("Q22").Value = [double]"0.8775106357116667"
$ws.Range("R22").Value = [double]"7.897595721405"
$ws.Range("S22").Value = [double]"2.458205511752763E-05"
$ws.Range("T22").Value = [double]"2.458205511752762E-05"
$ws.Range("G23").Value = [double]"0.178139"
$ws.Range("H23").Value = [double]"0.534417"
$ws.Range("I23").Value = [double]"0.0004345431092138444"
$ws.Range("J23").Value = [double]"0.0004345431092138443"
$ws.Range("O23").Value = [double]"0.5464678959362861"
$ws.Range("P23").Value = [double]"0.5464678959362861"
$ws.Range("Q23").Value = [double]"8.476795797171"
$ws.Range("R23").Value = [double]"76.291162174539"
$ws.Range("S23").Value = [double]"0.0002374638585857013"
$ws.Range("T23").Value = [double]"0.0002374638585857013"
$ws.Range("G24").Value = [double]"0.178139"
$ws.Range("H24").Value = [double]"0.534417"
$ws.Range("I24").Value = [double]"0.0004345431092138444"
$ws.Range("J24").Value = [double]"0.0004345431092138443"
$ws.Range("M24").Value = [double]"13.80191933333334"
$ws.Range("N24").Value = [double]"41.40575800000001"
$ws.Range("O24").Value = [double]"0.1585007882996995"
$ws.Range("P24").Value = [double]"0.1585007882996994"
$ws.Range("Q24").Value = [double]"2.458660108120668"
$ws.Range("R24").Value = [double]"22.12794097308601"
$ws.Range("S24").Value = [double]"6.887542536059674E-05"
$ws.Range("T24").Value = [double]"6.887542536059671E-05"
$ws.Range("G25").Value = [double]"0.178139"
$ws.Range("H25").Value = [double]"0.534417"
$ws.Range("I25").Value = [double]"0.0004345431092138444"
$ws.Range("J25").Value = [double]"0.0004345431092138443"
$ws.Range("M25").Value = [double]"12.18796133333333"
$ws.Range("N25").Value = [double]"36.563884"
$ws.Range("O25").Value = [double]"0.1399661476381804"
$ws.Range("P25").Value = [double]"0.1399661476381803"
$ws.Range("Q25").Value = [double]"2.171151243958667"
$ws.Range("R25").Value = [double]"19.540361195628"
$ws.Range("S25").Value = [double]"6.082132497937888E-05"
$ws.Range("T25").Value = [double]"6.082132497937885E-05"
$ws.Range("G26").Value = [double]"0.178139"
$ws.Range("H26").Value = [double]"0.534417"
$ws.Range("I26").Value = [double]"0.0004345431092138444"
$ws.Range("J26").Value = [double]"0.0004345431092138443"
$ws.Range("M26").Value = [double]"8.576764333333333"
$ws.Range("N26").Value = [double]"25.730293"
$ws.Range("O26").Value = [double]"0.09849527990001386"
$ws.Range("P26").Value = [double]"0.09849527990001385"
$ws.Range("Q26").Value = [double]"1.527856221575667"
$ws.Range("R26").Value = [double]"13.750705994181"
$ws.Range("S26").Value = [double]"4.28004451706399E-05"
$ws.Range("T26").Value = [double]"4.280044517063988E-05"
